$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the rich-text formatted command strings (previously
# "EasyshellLib.getElement('MAIN_WINDOW').Exists()" with mixed run
# formatting) with the new plain-text command used by both rows that
# referenced it.
$ws.Range("B4").Value = "EasyShellTest().check_main_window(True)"
$ws.Range("B7").Value = "EasyShellTest().check_main_window(True)"
